$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the current contents of rows 2 and 3 (columns A, B, C, E, F) before
# touching anything. D is identical on both rows ("Advogado Exemplo ") so
# it is left untouched.
$row2A = $ws.Range("A2").Value2
$row2B = $ws.Range("B2").Value2
$row2C = $ws.Range("C2").Value2
$row2E = $ws.Range("E2").Value2
$row2F = $ws.Range("F2").Value2

$row3A = $ws.Range("A3").Value2
$row3B = $ws.Range("B3").Value2
$row3C = $ws.Range("C3").Value2
$row3E = $ws.Range("E3").Value2
$row3F = $ws.Range("F3").Value2

# Helper scratch cell used to stage every write. Writing straight into the
# target cell via ".Value = <text>" lets Excel's normal type inference
# kick in, which would silently turn text that merely looks like a date
# (e.g. "7/8/2024") into a real date serial number. Building the literal
# via a formula that evaluates to a text string, then copying just the
# computed VALUE (PasteSpecial xlPasteValues = -4163) into the target
# cell, preserves the original cell style and stores the text verbatim
# (matching the inlineStr cells produced by the original export).
$scratch = $ws.Range("Z1")

function Set-LiteralValue($cell, $value) {
    if ($value -is [string]) {
        $escaped = $value.Replace('"', '""')
        $scratch.Formula = '="' + $escaped + '"'
    } else {
        $scratch.Formula = "=" + $value
    }
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

Set-LiteralValue $ws.Range("A2") $row3A
Set-LiteralValue $ws.Range("B2") $row3B
Set-LiteralValue $ws.Range("C2") $row3C
Set-LiteralValue $ws.Range("E2") $row3E
Set-LiteralValue $ws.Range("F2") $row3F

Set-LiteralValue $ws.Range("A3") $row2A
Set-LiteralValue $ws.Range("B3") $row2B
Set-LiteralValue $ws.Range("C3") $row2C
Set-LiteralValue $ws.Range("E3") $row2E
Set-LiteralValue $ws.Range("F3") $row2F

$scratch.Clear()
